$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1. Rename "LP 4" -> "LP 5" on the slide that currently holds it.
# ------------------------------------------------------------------
$lpSlide = $p.Slides.Item(29)
$lpSlide.Shapes.Item(1).TextFrame.TextRange.Text = "LP 5"

# ------------------------------------------------------------------
# 2. Duplicate the "Alphabet of Lines" / "Drawing #8" slide (slide 31)
#    to create a new Unit-4 divider slide, move it to slide 30, and
#    retitle it.
# ------------------------------------------------------------------
$srcDivider = $p.Slides.Item(31)
$dividerRange = $srcDivider.Duplicate()
$divider = $dividerRange.Item(1)
$divider.MoveTo(30)

$divider.Shapes.Item(1).TextFrame.TextRange.Text = "Unit 4 Alphabet of Lines"
$divider.Shapes.Item(2).TextFrame.TextRange.Text = ""

# ------------------------------------------------------------------
# 3. Duplicate the "Google classroom code" slide (slide 2) and move
#    the copy to slide 31 (right after the new divider slide).
# ------------------------------------------------------------------
$srcCode = $p.Slides.Item(2)
$codeRange = $srcCode.Duplicate()
$code = $codeRange.Item(1)
$code.MoveTo(31)
